$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.393.27"
$cell.Style = $origStyle
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.641.58"
$cell.Style = $origStyle
$ws.Range("E3").Value = "  -1.53%  "

# Row 4
$cell = $ws.Range("D4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "211.83"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -1.56%  "

# Row 6
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.535"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +4.10%  "

# Row 7
$cell = $ws.Range("D7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "23.19"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  -1.69%  "

# Row 9
$ws.Range("E9").Value = "  -2.31%  "

# Row 10
$ws.Range("E10").Value = "  -2.03%  "

# Row 11
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0893"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  +1.48%  "

# Row 12
$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.874.65"
$cell.Style = $origStyle
$ws.Range("E12").Value = "  -1.46%  "

# Row 13
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.641.79"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  -1.44%  "

# Row 14
$ws.Range("E14").Value = "  -3.14%  "

# Row 15
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.559"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +0.76%  "

# Row 16
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "64.23"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  -3.08%  "

# Row 17
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "27.390.36"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  -0.70%  "

# Row 18
$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "227.97"
$cell.Style = $origStyle
$ws.Range("E18").Value = "  -9.32%  "

# Row 19
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0718"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  -2.05%  "

# Row 20
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.45"
$cell.Style = $origStyle
$ws.Range("E20").Value = "  -1.07%  "

# Row 21
$ws.Range("E21").Value = "  -0.05%  "

# Row 22
$ws.Range("E22").Value = "  -4.26%  "

# Row 23
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.15"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  -1.58%  "

# Row 24
$ws.Range("E24").Value = "  +0.32%  "

# Row 25
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "147.63"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  +0.67%  "

# Row 26
$ws.Range("E26").Value = "  +2.23%  "

# Row 27
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.93"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -2.98%  "

# Row 28
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  -0.06%  "

# Row 29
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.50"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -6.50%  "

# Row 30
$ws.Range("E30").Value = "  -4.75%  "

# Row 31
$ws.Range("E31").Value = "  -4.56%  "

# Row 32
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.26"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -2.79%  "

# Row 33
$ws.Range("E33").Value = "  -0.65%  "

# Row 34
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.396.48"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  -5.45%  "

# Row 35
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.55"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  -1.21%  "

# Row 36
$ws.Range("E36").Value = "  -0.39%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.560"
$cell.Style = $origStyle
$ws.Range("E37").Value = "  -3.16%  "

# Row 38
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.878"
$cell.Style = $origStyle
$ws.Range("E38").Value = "  -6.93%  "

# Row 39
$ws.Range("E39").Value = "  -3.16%  "

# Row 40
$ws.Range("E40").Value = "  +0.42%  "

# Row 41
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.46"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  -2.28%  "

# Row 43
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.47"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +0.77%  "

# Row 44
$ws.Range("E44").Value = "  +0.00%  "

# Row 45
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.787"
$cell.Style = $origStyle
$ws.Range("E45").Value = "  -0.66%  "

# Row 46
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "64.21"
$cell.Style = $origStyle
$ws.Range("E46").Value = "  -7.97%  "

# Row 47
$cell = $ws.Range("D47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.785.78"
$cell.Style = $origStyle
$ws.Range("E47").Value = "  -1.33%  "

# Row 48
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.64"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  -3.66%  "

# Row 49
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "87.18"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  -2.59%  "

# Row 50
$ws.Range("E50").Value = "  -4.32%  "

# Row 51
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0982"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  -3.54%  "
